$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Formula = "=`"29.375.28`""
$ws.Range("D2").Copy()
$ws.Range("D2").PasteSpecial(-4163)
$ws.Range("E2").Value = "  -0.09%  "
$ws.Range("D3").Formula = "=`"1.849.49`""
$ws.Range("D3").Copy()
$ws.Range("D3").PasteSpecial(-4163)
$ws.Range("E3").Value = "  -0.06%  "
$ws.Range("D4").Formula = "=`"1.001`""
$ws.Range("D4").Copy()
$ws.Range("D4").PasteSpecial(-4163)
$ws.Range("E4").Value = "  +0.14%  "
$ws.Range("D5").Formula = "=`"240.47`""
$ws.Range("D5").Copy()
$ws.Range("D5").PasteSpecial(-4163)
$ws.Range("E5").Value = "  -0.10%  "
$ws.Range("D6").Formula = "=`"0.6279`""
$ws.Range("D6").Copy()
$ws.Range("D6").PasteSpecial(-4163)
$ws.Range("E6").Value = "  -0.21%  "
$ws.Range("D7").Formula = "=`"1.000`""
$ws.Range("D7").Copy()
$ws.Range("D7").PasteSpecial(-4163)
$ws.Range("E7").Value = "  +0.02%  "
$ws.Range("D8").Formula = "=`"0.07608`""
$ws.Range("D8").Copy()
$ws.Range("D8").PasteSpecial(-4163)
$ws.Range("E8").Value = "  -0.87%  "
$ws.Range("D9").Formula = "=`"0.2912`""
$ws.Range("D9").Copy()
$ws.Range("D9").PasteSpecial(-4163)
$ws.Range("E9").Value = "  -0.76%  "
$ws.Range("D10").Formula = "=`"24.58`""
$ws.Range("D10").Copy()
$ws.Range("D10").PasteSpecial(-4163)
$ws.Range("E10").Value = "  +0.12%  "
$ws.Range("D11").Formula = "=`"0.07753`""
$ws.Range("D11").Copy()
$ws.Range("D11").PasteSpecial(-4163)
$ws.Range("E11").Value = "  +0.07%  "
$ws.Range("D12").Formula = "=`"5.020`""
$ws.Range("D12").Copy()
$ws.Range("D12").PasteSpecial(-4163)
$ws.Range("E12").Value = "  -0.02%  "
$ws.Range("D13").Formula = "=`"0.6806`""
$ws.Range("D13").Copy()
$ws.Range("D13").PasteSpecial(-4163)
$ws.Range("E13").Value = "  -0.02%  "
$ws.Range("D14").Formula = "=`"0.00001049`""
$ws.Range("D14").Copy()
$ws.Range("D14").PasteSpecial(-4163)
$ws.Range("E14").Value = "  -5.76%  "
$ws.Range("D15").Formula = "=`"83.21`""
$ws.Range("D15").Copy()
$ws.Range("D15").PasteSpecial(-4163)
$ws.Range("E15").Value = "  -0.45%  "
$ws.Range("D16").Formula = "=`"6.128`""
$ws.Range("D16").Copy()
$ws.Range("D16").PasteSpecial(-4163)
$ws.Range("E16").Value = "  -0.29%  "
$ws.Range("D17").Formula = "=`"29.392.63`""
$ws.Range("D17").Copy()
$ws.Range("D17").PasteSpecial(-4163)
$ws.Range("E17").Value = "  -0.24%  "
$ws.Range("D18").Formula = "=`"229.33`""
$ws.Range("D18").Copy()
$ws.Range("D18").PasteSpecial(-4163)
$ws.Range("E18").Value = "  +0.10%  "
$ws.Range("D19").Formula = "=`"12.33`""
$ws.Range("D19").Copy()
$ws.Range("D19").PasteSpecial(-4163)
$ws.Range("E19").Value = "  -1.08%  "
$ws.Range("D20").Formula = "=`"1.001`""
$ws.Range("D20").Copy()
$ws.Range("D20").PasteSpecial(-4163)
$ws.Range("E20").Value = "  +0.04%  "
$ws.Range("D21").Formula = "=`"7.472`""
$ws.Range("D21").Copy()
$ws.Range("D21").PasteSpecial(-4163)
$ws.Range("E21").Value = "  +0.42%  "
$ws.Range("D22").Formula = "=`"1.000`""
$ws.Range("D22").Copy()
$ws.Range("D22").PasteSpecial(-4163)
$ws.Range("E22").Value = "  +0.01%  "
$ws.Range("D23").Formula = "=`"158.76`""
$ws.Range("D23").Copy()
$ws.Range("D23").PasteSpecial(-4163)
$ws.Range("E23").Value = "  +1.20%  "
$ws.Range("D24").Formula = "=`"0.1389`""
$ws.Range("D24").Copy()
$ws.Range("D24").PasteSpecial(-4163)
$ws.Range("E24").Value = "  +0.32%  "
$ws.Range("D25").Formula = "=`"8.439`""
$ws.Range("D25").Copy()
$ws.Range("D25").PasteSpecial(-4163)
$ws.Range("E25").Value = "  +0.67%  "
$ws.Range("D26").Formula = "=`"17.69`""
$ws.Range("D26").Copy()
$ws.Range("D26").PasteSpecial(-4163)
$ws.Range("E26").Value = "  +0.10%  "
$ws.Range("E27").Value = "  +10.37%  "
$ws.Range("D28").Formula = "=`"1.475`""
$ws.Range("D28").Copy()
$ws.Range("D28").PasteSpecial(-4163)
$ws.Range("E28").Value = "  +0.49%  "
$ws.Range("D29").Formula = "=`"0.05601`""
$ws.Range("D29").Copy()
$ws.Range("D29").PasteSpecial(-4163)
$ws.Range("E29").Value = "  -2.13%  "
$ws.Range("D30").Formula = "=`"4.110`""
$ws.Range("D30").Copy()
$ws.Range("D30").PasteSpecial(-4163)
$ws.Range("E30").Value = "  -0.37%  "
$ws.Range("E31").Value = "  +0.47%  "
$ws.Range("D32").Formula = "=`"1.833`""
$ws.Range("D32").Copy()
$ws.Range("D32").PasteSpecial(-4163)
$ws.Range("E32").Value = "  -0.81%  "
$ws.Range("D33").Formula = "=`"1.156`""
$ws.Range("D33").Copy()
$ws.Range("D33").PasteSpecial(-4163)
$ws.Range("E33").Value = "  -0.53%  "
$ws.Range("D34").Formula = "=`"0.7021`""
$ws.Range("D34").Copy()
$ws.Range("D34").PasteSpecial(-4163)
$ws.Range("E34").Value = "  -0.88%  "
$ws.Range("D35").Formula = "=`"2.585`""
$ws.Range("D35").Copy()
$ws.Range("D35").PasteSpecial(-4163)
$ws.Range("E35").Value = "  +0.10%  "
$ws.Range("D36").Formula = "=`"1.236.71`""
$ws.Range("D36").Copy()
$ws.Range("D36").PasteSpecial(-4163)
$ws.Range("E36").Value = "  +1.64%  "
$ws.Range("D37").Formula = "=`"0.01804`""
$ws.Range("D37").Copy()
$ws.Range("D37").PasteSpecial(-4163)
$ws.Range("E37").Value = "  +0.69%  "
$ws.Range("D38").Formula = "=`"2.731`""
$ws.Range("D38").Copy()
$ws.Range("D38").PasteSpecial(-4163)
$ws.Range("E38").Value = "  -1.63%  "
$ws.Range("D39").Formula = "=`"6.426`""
$ws.Range("D39").Copy()
$ws.Range("D39").PasteSpecial(-4163)
$ws.Range("E39").Value = "  -1.16%  "
$ws.Range("D40").Formula = "=`"0.9032`""
$ws.Range("D40").Copy()
$ws.Range("D40").PasteSpecial(-4163)
$ws.Range("E40").Value = "  -0.46%  "
$ws.Range("D41").Formula = "=`"0.9997`""
$ws.Range("D41").Copy()
$ws.Range("D41").PasteSpecial(-4163)
$ws.Range("E41").Value = "  -0.05%  "
$ws.Range("D42").Formula = "=`"101.58`""
$ws.Range("D42").Copy()
$ws.Range("D42").PasteSpecial(-4163)
$ws.Range("E42").Value = "  -0.24%  "
$ws.Range("D43").Formula = "=`"65.55`""
$ws.Range("D43").Copy()
$ws.Range("D43").PasteSpecial(-4163)
$ws.Range("E43").Value = "  -1.21%  "
$ws.Range("D44").Formula = "=`"7.182`""
$ws.Range("D44").Copy()
$ws.Range("D44").PasteSpecial(-4163)
$ws.Range("E44").Value = "  +0.79%  "
$ws.Range("B45").Value = "BabyDogeCoin"
$ws.Range("C45").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D45").Formula = "=`"0.00000000118`""
$ws.Range("D45").Copy()
$ws.Range("D45").PasteSpecial(-4163)
$ws.Range("E45").Value = "  -3.03%  "
$ws.Range("B46").Value = "TheSandbox"
$ws.Range("C46").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D46").Formula = "=`"0.3997`""
$ws.Range("D46").Copy()
$ws.Range("D46").PasteSpecial(-4163)
$ws.Range("E46").Value = "  -0.42%  "
$ws.Range("B47").Value = "Algorand"
$ws.Range("C47").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D47").Formula = "=`"0.1157`""
$ws.Range("D47").Copy()
$ws.Range("D47").PasteSpecial(-4163)
$ws.Range("E47").Value = "  +2.32%  "
$ws.Range("D48").Formula = "=`"1.686`""
$ws.Range("D48").Copy()
$ws.Range("D48").PasteSpecial(-4163)
$ws.Range("E48").Value = "  +0.22%  "
$ws.Range("B49").Value = "EnergySwap"
$ws.Range("C49").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D49").Formula = "=`"8.991`""
$ws.Range("D49").Copy()
$ws.Range("D49").PasteSpecial(-4163)
$ws.Range("E49").Value = "  +0.17%  "
$ws.Range("D50").Formula = "=`"0.05702`""
$ws.Range("D50").Copy()
$ws.Range("D50").PasteSpecial(-4163)
$ws.Range("E50").Value = "  -0.13%  "
$ws.Range("D51").Formula = "=`"0.4632`""
$ws.Range("D51").Copy()
$ws.Range("D51").PasteSpecial(-4163)
$ws.Range("E51").Value = "  +0.02%  "

$excel.CutCopyMode = $false
